# Auto-generated Excel COM-interop script
# Applies cached-price / profit-column updates per the scraped OOXML diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 266.72223
$ws.Range("I5").Value = 219.57143
$ws.Range("K5").Value = 219.57143
$ws.Range("M5").Value = -104.57143

$ws.Range("H19").Value = 1738.0769
$ws.Range("I19").Value = 861.75
$ws.Range("K19").Value = 861.75
$ws.Range("M19").Value = -686.75

$ws.Range("H43").Value = 14825
$ws.Range("J43").Value = 8333.333000000001
$ws.Range("L43").Value = 8333.333000000001
$ws.Range("N43").Value = -8471.333000000001

$ws.Range("H100").Value = 7050.2666
$ws.Range("I100").Value = 2089.7144
$ws.Range("J100").Value = 9290.517
$ws.Range("K100").Value = 2089.7144
$ws.Range("L100").Value = 9290.517
$ws.Range("M100").Value = -1548.7144
$ws.Range("N100").Value = -10372.517

$ws.Range("H112").Value = 2669.6
$ws.Range("J112").Value = 2669.6
$ws.Range("L112").Value = 8008.799999999999
$ws.Range("N112").Value = -10224.8

$ws.Range("H138").Value = 2968.6858
$ws.Range("J138").Value = 3863.8235
$ws.Range("L138").Value = 11591.4705
$ws.Range("N138").Value = -21871.4705

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 902.5714
$ws.Range("I8").Value = 553
$ws.Range("J8").Value = 1042.4
$ws.Range("K8").Value = 553
$ws.Range("L8").Value = 1042.4
$ws.Range("M8").Value = -409
$ws.Range("N8").Value = -1330.4

$ws.Range("H18").Value = 5580
$ws.Range("J18").Value = 5580
$ws.Range("L18").Value = 5580
$ws.Range("N18").Value = -6224

$ws.Range("H32").Value = 13473.1875
$ws.Range("I32").Value = 13517.152
$ws.Range("K32").Value = 13517.152
$ws.Range("M32").Value = -13230.152

$ws.Range("H61").Value = 4041.6038
$ws.Range("I61").Value = 3368.6382
$ws.Range("K61").Value = 3368.6382
$ws.Range("M61").Value = -3156.6382

$ws.Range("H74").Value = 406439.28
$ws.Range("I74").Value = 459999.2
$ws.Range("J74").Value = 13666.667
$ws.Range("K74").Value = 459999.2
$ws.Range("L74").Value = 13666.667
$ws.Range("M74").Value = -459125.2
$ws.Range("N74").Value = -15414.667

$ws.Range("H77").Value = 406439.28
$ws.Range("I77").Value = 459999.2
$ws.Range("J77").Value = 13666.667
$ws.Range("K77").Value = 2299996
$ws.Range("L77").Value = 68333.33499999999
$ws.Range("M77").Value = -2295628
$ws.Range("N77").Value = -77069.33499999999

$ws.Range("H97").Value = 1159105.1
$ws.Range("I97").Value = 1278741.5
$ws.Range("K97").Value = 1278741.5
$ws.Range("M97").Value = -1278245.5

$ws.Range("H132").Value = 7197.8945
$ws.Range("I132").Value = 5774.7856
$ws.Range("K132").Value = 17324.3568
$ws.Range("M132").Value = -14794.3568

$ws.Range("H136").Value = 4041.6038
$ws.Range("I136").Value = 3368.6382
$ws.Range("K136").Value = 10105.9146
$ws.Range("M136").Value = -7555.9146

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2976.12
$ws.Range("I20").Value = 2770.6428
$ws.Range("K20").Value = 2770.6428
$ws.Range("M20").Value = -2523.6428

$ws.Range("H134").Value = 6733.3335
$ws.Range("I134").Value = 4284.625
$ws.Range("K134").Value = 12853.875
$ws.Range("M134").Value = -10318.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 2253
$ws.Range("I17").Value = 2253
$ws.Range("K17").Value = 2253
$ws.Range("M17").Value = -2079

$ws.Range("H22").Value = 1322.1765
$ws.Range("I22").Value = 652.7273
$ws.Range("J22").Value = 2549.5
$ws.Range("K22").Value = 652.7273
$ws.Range("L22").Value = 2549.5
$ws.Range("M22").Value = -302.7273
$ws.Range("N22").Value = -3249.5

$ws.Range("H25").Value = 24300
$ws.Range("I25").Value = 36929
$ws.Range("J25").Value = 11671
$ws.Range("K25").Value = 36929
$ws.Range("L25").Value = 11671
$ws.Range("M25").Value = -36755
$ws.Range("N25").Value = -12019

$ws.Range("H87").Value = 60000
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 60000
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 60000
$ws.Range("N87").Value = -62372
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value = 60000
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 60000
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 180000
$ws.Range("N90").Value = -191856
$ws.Range("M90").ClearContents()

$ws.Range("H132").Value = 3311.1177
$ws.Range("I132").Value = 2486.8125
$ws.Range("K132").Value = 7460.4375
$ws.Range("M132").Value = -4930.4375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 35176824
$ws.Range("I4").Value = 48557310
$ws.Range("K4").Value = 145671930
$ws.Range("M4").Value = -145671818

$ws.Range("H5").Value = 1290.4546
$ws.Range("I5").Value = 800
$ws.Range("K5").Value = 2400
$ws.Range("M5").Value = -2288

$ws.Range("H13").Value = 268.2
$ws.Range("I13").Value = 179.66667
$ws.Range("J13").Value = 401
$ws.Range("K13").Value = 539.00001
$ws.Range("L13").Value = 1203
$ws.Range("M13").Value = -371.00001
$ws.Range("N13").Value = -1539

$ws.Range("H20").Value = 2691.1667
$ws.Range("J20").Value = 3349.25
$ws.Range("L20").Value = 10047.75
$ws.Range("N20").Value = -10501.75

$ws.Range("H80").Value = 3650
$ws.Range("J80").Value = 3780
$ws.Range("L80").Value = 11340
$ws.Range("N80").Value = -13212

$ws.Range("H83").Value = 3650
$ws.Range("J83").Value = 3780
$ws.Range("L83").Value = 34020
$ws.Range("N83").Value = -43380

$ws.Range("H92").Value = 1691.6666
$ws.Range("I92").Value = 1725
$ws.Range("J92").Value = 1650
$ws.Range("K92").Value = 5175
$ws.Range("L92").Value = 4950
$ws.Range("M92").Value = -3927
$ws.Range("N92").Value = -7446

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H113").Value = 1593.4
$ws.Range("I113").Value = 1554.5
$ws.Range("J113").Value = 1612.85
$ws.Range("K113").Value = 4663.5
$ws.Range("L113").Value = 4838.549999999999
$ws.Range("M113").Value = -2493.5
$ws.Range("N113").Value = -9178.549999999999

$ws.Range("H135").Value = 1290.4546
$ws.Range("I135").Value = 800
$ws.Range("K135").Value = 7200
$ws.Range("M135").Value = -4665

$ws.Range("H137").Value = 3013.5386
$ws.Range("J137").Value = 3116.0908
$ws.Range("L137").Value = 9348.2724
$ws.Range("N137").Value = -19548.2724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 579.0476
$ws.Range("J2").Value = 155.33333
$ws.Range("L2").Value = 155.33333
$ws.Range("N2").Value = -381.33333

$ws.Range("H80").Value = 5976
$ws.Range("I80").Value = 4973.625
$ws.Range("J80").Value = 7121.5713
$ws.Range("K80").Value = 4973.625
$ws.Range("L80").Value = 7121.5713
$ws.Range("M80").Value = -3975.625
$ws.Range("N80").Value = -9117.5713

$ws.Range("H83").Value = 5976
$ws.Range("I83").Value = 4973.625
$ws.Range("J83").Value = 7121.5713
$ws.Range("K83").Value = 24868.125
$ws.Range("L83").Value = 35607.85649999999
$ws.Range("M83").Value = -19876.125
$ws.Range("N83").Value = -45591.85649999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3833
$ws.Range("J7").Value = 3750
$ws.Range("L7").Value = 3750
$ws.Range("N7").Value = -3974

$ws.Range("H35").Value = 963.5
$ws.Range("I35").Value = 583
$ws.Range("J35").Value = 2105
$ws.Range("K35").Value = 583
$ws.Range("L35").Value = 2105
$ws.Range("M35").Value = -247
$ws.Range("N35").Value = -2777

$ws.Range("H68").Value = 3827.7778
$ws.Range("I68").Value = 3200
$ws.Range("J68").Value = 5083.3335
$ws.Range("K68").Value = 3200
$ws.Range("L68").Value = 5083.3335
$ws.Range("M68").Value = -2451
$ws.Range("N68").Value = -6581.3335

$ws.Range("H71").Value = 3827.7778
$ws.Range("I71").Value = 3200
$ws.Range("J71").Value = 5083.3335
$ws.Range("K71").Value = 16000
$ws.Range("L71").Value = 25416.6675
$ws.Range("M71").Value = -12256
$ws.Range("N71").Value = -32904.6675

$ws.Range("H93").Value = 1592
$ws.Range("I93").Value = 1724.7778
$ws.Range("J93").Value = 397
$ws.Range("K93").Value = 1724.7778
$ws.Range("L93").Value = 397
$ws.Range("M93").Value = -476.7778000000001
$ws.Range("N93").Value = -2893

$ws.Range("H126").Value = 3833
$ws.Range("J126").Value = 3750
$ws.Range("L126").Value = 11250
$ws.Range("N126").Value = -16190

$ws.Range("H132").Value = 7111.1
$ws.Range("I132").Value = 5058.1875
$ws.Range("J132").Value = 9457.286
$ws.Range("K132").Value = 15174.5625
$ws.Range("L132").Value = 28371.858
$ws.Range("M132").Value = -12644.5625
$ws.Range("N132").Value = -33431.858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 31883
$ws.Range("J92").Value = 31883
$ws.Range("L92").Value = 31883
$ws.Range("N92").Value = -36875

$ws.Range("H132").Value = 5201.517
$ws.Range("I132").Value = 4065.75
$ws.Range("J132").Value = 8771.071
$ws.Range("K132").Value = 12197.25
$ws.Range("L132").Value = 26313.213
$ws.Range("M132").Value = -9667.25
$ws.Range("N132").Value = -31373.213
